$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.080.57"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "'3.741.40"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'600.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.01%  "
$ws.Range("D6").Value = "'167.18"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").Value = "'3.740.45"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +1.39%  "
$ws.Range("D10").Value = "'0.169"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.67%  "
$ws.Range("E11").Value = "  +0.57%  "
$ws.Range("D12").Value = "'0.460"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.32%  "
$ws.Range("D13").Value = "'37.97"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.13%  "
$ws.Range("E14").Value = "  +1.46%  "
$ws.Range("D15").Value = "'4.367.12"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.16%  "
$ws.Range("D16").Value = "'3.738.50"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.05%  "
$ws.Range("D17").Value = "'69.044.45"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.31%  "
$ws.Range("E18").Value = "  +1.26%  "
$ws.Range("D19").Value = "'17.32"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.29%  "
$ws.Range("E20").Value = "  -1.48%  "
$ws.Range("D21").Value = "'11.18"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +11.44%  "
$ws.Range("D22").Value = "'492.02"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.96%  "
$ws.Range("D23").Value = "'0.727"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.51%  "
$ws.Range("E24").Value = "  +8.32%  "
$ws.Range("D25").Value = "'84.94"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.28%  "
$ws.Range("D26").Value = "'2.29"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.69%  "
$ws.Range("D27").Value = "'12.21"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.09%  "
$ws.Range("D28").Value = "'10.05"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.66%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").Value = "'2.97"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.10%  "
$ws.Range("D31").Value = "'8.17"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.59%  "
$ws.Range("E32").Value = "  -0.19%  "
$ws.Range("D33").Value = "'31.41"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.87%  "
$ws.Range("D34").Value = "'3.887.51"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.14%  "
$ws.Range("D35").Value = "'3.674.63"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.46%  "
$ws.Range("E36").Value = "  -0.33%  "
$ws.Range("D37").Value = "'1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("D38").Value = "'1.01"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Value = "'5.94"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.83%  "
$ws.Range("E40").Value = "  +5.05%  "
$ws.Range("E41").Value = "  -0.17%  "
$ws.Range("E42").Value = "  +6.19%  "
$ws.Range("D44").Value = "'1.98"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.52%  "
$ws.Range("D45").Value = "'422.25"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.71%  "
$ws.Range("D46").Value = "'8.44"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.44%  "
$ws.Range("D48").Value = "'40.01"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.50%  "
$ws.Range("D49").Value = "'141.85"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.35%  "
$ws.Range("D50").Value = "'2.780.46"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.42%  "
$ws.Range("D51").Value = "'0.0352"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.07%  "
